$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44995
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 5500
$ws.Range("P2").Value = 5750
$ws.Range("S2").Value = 2875

# Row 3
$ws.Range("D3").Value = 45008
$ws.Range("L3").Value = "Especial"
$ws.Range("M3").Value = 60
$ws.Range("N3").Value = 7000
$ws.Range("O3").Value = 7000
$ws.Range("P3").Value = 7000
$ws.Range("S3").Value = 3500

# Row 4
$ws.Range("L4").Value = "Primera"
$ws.Range("N4").Value = 6000
$ws.Range("O4").Value = 6000
$ws.Range("P4").Value = 6000
$ws.Range("S4").Value = 3000

# Row 5
$ws.Range("D5").Value = 44991
$ws.Range("M5").Value = 50
